$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 00:35"

# --- Estados Unidos (row 4): refresh daily figures ---
$ws.Range("B4").Value = 1569033
$ws.Range("C4").Value = 18739
$ws.Range("D4").Value = 362591
$ws.Range("E4").Value = 1112997
$ws.Range("G4").Value = 1464
$ws.Range("H4").Value = 93445

# --- Brasil (row 7): refresh daily figures ---
$ws.Range("B7").Value = 271628
$ws.Range("C7").Value = 16260
$ws.Range("D7").Value = 106794
$ws.Range("E7").Value = 146863
$ws.Range("G7").Value = 1118
$ws.Range("H7").Value = 17971

# --- Insert a new row for Colombia right after Rumania (row 38), before Kuwait ---
$ws.Rows("39").Insert()
$ws.Range("A39").Value = "Colombia"
$ws.Range("B39").Value = 16935
$ws.Range("C39").Value = 640
$ws.Range("D39").Value = 4050
$ws.Range("E39").Value = 12272
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 21
$ws.Range("H39").Value = 613

# The old Colombia entry (now pushed down to row 44, after Japon) is removed
# since Colombia now appears earlier in the list.
$ws.Rows("44").Delete()

# --- Noruega (row 53): refresh daily figures ---
$ws.Range("B53").Value = 8267
$ws.Range("C53").Value = 10
$ws.Range("E53").Value = 8002

# --- Guinea-Bisau (row 104): refresh daily figures ---
$ws.Range("B104").Value = 1038
$ws.Range("C104").Value = 6
$ws.Range("D104").Value = 42
$ws.Range("G104").Value = 2
$ws.Range("H104").Value = 6

# --- Niger (row 111): refresh daily figures ---
$ws.Range("B111").Value = 914
$ws.Range("C111").Value = 5
$ws.Range("D111").Value = 734
$ws.Range("E111").Value = 125
